$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D8").Value = 247000
$ws.Range("E8").Value = 218700
$ws.Range("F8").Value = 155400
$ws.Range("G8").Value = 195000
$ws.Range("H8").Value = 239100
$ws.Range("I8").Value = 272300
$ws.Range("J8").Value = 351700

$ws.Range("D9").Value = 204500
$ws.Range("E9").Value = 180300
$ws.Range("F9").Value = 124800
$ws.Range("G9").Value = 176000
$ws.Range("H9").Value = 202000
$ws.Range("I9").Value = 201000
$ws.Range("J9").Value = 500700

$ws.Range("D10").Value = 42600
$ws.Range("E10").Value = 38400
$ws.Range("F10").Value = 30600
$ws.Range("G10").Value = 19000
$ws.Range("H10").Value = 37000
$ws.Range("I10").Value = 71300
$ws.Range("J10").Value = -148900

$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 22600
$ws.Range("F14").Value = 1000
$ws.Range("G14").Value = 31600
$ws.Range("H14").Value = 45300
$ws.Range("I14").Value = 1700
$ws.Range("J14").Value = 0

$ws.Range("D15").Value = 22300
$ws.Range("E15").Value = 21100
$ws.Range("F15").Value = 22600
$ws.Range("G15").Value = 31600
$ws.Range("H15").Value = 37200
$ws.Range("I15").Value = 39600
$ws.Range("J15").Value = 41300

$ws.Range("D17").Value = 254700
$ws.Range("E17").Value = 235400
$ws.Range("F17").Value = 164000
$ws.Range("G17").Value = 253800
$ws.Range("H17").Value = 303200
$ws.Range("I17").Value = 251700
$ws.Range("J17").Value = 310400

$ws.Range("D18").Value = -7600
$ws.Range("E18").Value = -16700
$ws.Range("F18").Value = -8600
$ws.Range("G18").Value = -58700
$ws.Range("H18").Value = -64100
$ws.Range("I18").Value = 20500
$ws.Range("J18").Value = 41400

$ws.Range("D20").Value = 100
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = -200
$ws.Range("G20").Value = 400
$ws.Range("H20").Value = -11400
$ws.Range("I20").Value = -4500
$ws.Range("J20").Value = 5800

$ws.Range("D21").Value = 14100
$ws.Range("E21").Value = 7600
$ws.Range("F21").Value = 24400
$ws.Range("G21").Value = -21100
$ws.Range("H21").Value = -33000
$ws.Range("I21").Value = 57400
$ws.Range("J21").Value = "NA"

$ws.Range("D22").Value = 8500
$ws.Range("E22").Value = 4000
$ws.Range("F22").Value = 6500
$ws.Range("G22").Value = 5000
$ws.Range("H22").Value = 5700
$ws.Range("I22").Value = 7900
$ws.Range("J22").Value = 7800

$ws.Range("D23").Value = -16000
$ws.Range("E23").Value = -20600
$ws.Range("F23").Value = -15300
$ws.Range("G23").Value = -63300
$ws.Range("H23").Value = -81300
$ws.Range("I23").Value = 8100
$ws.Range("J23").Value = 39400

$ws.Range("D26").Value = "NA"
$ws.Range("E26").Value = -20600
$ws.Range("F26").Value = -15300
$ws.Range("G26").Value = -63300
$ws.Range("H26").Value = -81300
$ws.Range("I26").Value = 8100
$ws.Range("J26").Value = 39400

$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = -24500
$ws.Range("F27").Value = 6300
$ws.Range("G27").Value = -39200
$ws.Range("H27").Value = -102800
$ws.Range("I27").Value = 3900
$ws.Range("J27").Value = 21300

$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 1800
$ws.Range("F29").Value = -115500
$ws.Range("G29").Value = 8100
$ws.Range("H29").Value = 130300
$ws.Range("I29").Value = 1300
$ws.Range("J29").Value = 100

$ws.Range("D32").Value = -100
$ws.Range("E32").Value = -100
$ws.Range("F32").Value = 200
$ws.Range("G32").Value = -400
$ws.Range("H32").Value = 11400
$ws.Range("I32").Value = 4500
$ws.Range("J32").Value = -5800

$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = -22700
$ws.Range("F33").Value = -109200
$ws.Range("G33").Value = -31100
$ws.Range("H33").Value = 27600
$ws.Range("I33").Value = 5200
$ws.Range("J33").Value = 21400

$ws.Range("D35").Value = "NA"
$ws.Range("E35").Value = -22700
$ws.Range("F35").Value = -109200
$ws.Range("G35").Value = -31100
$ws.Range("H35").Value = 27600
$ws.Range("I35").Value = 5200
$ws.Range("J35").Value = 21400

$ws.Range("D81").Value = "NA"
$ws.Range("E81").Value = -22700
$ws.Range("F81").Value = -109200
$ws.Range("G81").Value = -31100
$ws.Range("H81").Value = 27600
$ws.Range("I81").Value = 5200
$ws.Range("J81").Value = 21400

$ws.Range("J83").Value = "NA"

$ws.Range("J94").Value = "NA"

$ws.Range("J100").Value = "NA"

$ws.Range("J101").Value = "NA"
